# Auto-generated Excel COM script applying scheduled market-price refresh
# to the Sagittarius_Profits leve-profit tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 5697.7085
$ws.Range("J17").Value = 6304.524
$ws.Range("L17").Value = 18913.572
$ws.Range("N17").Value = -19249.572

# Row 33
$ws.Range("H33").Value = 1027.6666
$ws.Range("I33").Value = 1027.6666
$ws.Range("K33").Value = 1027.6666
$ws.Range("M33").Value = -798.6666

# Row 93
$ws.Range("H93").Value = 53331.668
$ws.Range("J93").Value = 53331.668
$ws.Range("L93").Value = 53331.668
$ws.Range("N93").Value = -58323.668

# Row 107
$ws.Range("H107").Value = 209
$ws.Range("I107").Value = 238
$ws.Range("K107").Value = 238
$ws.Range("M107").Value = 1682

# Row 127
$ws.Range("H127").Value = 3035.7144
$ws.Range("I127").Value = 3208.6667
$ws.Range("K127").Value = 9626.000100000001
$ws.Range("M127").Value = -4666.000100000001

# Row 138
$ws.Range("H138").Value = 2197.6553
$ws.Range("J138").Value = 3184.6667
$ws.Range("L138").Value = 9554.000100000001
$ws.Range("N138").Value = -19834.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4336.4443
$ws.Range("I32").Value = 4521.8823
$ws.Range("K32").Value = 4521.8823
$ws.Range("M32").Value = -4234.8823

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 61
$ws.Range("H61").Value = 2818.8667
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 97
$ws.Range("H97").Value = 755.625
$ws.Range("I97").Value = 755.625
$ws.Range("K97").Value = 755.625
$ws.Range("M97").Value = -259.625

# Row 102
$ws.Range("H102").Value = 1126
$ws.Range("I102").Value = 1126
$ws.Range("K102").Value = 1126
$ws.Range("M102").Value = 496

# Row 122
$ws.Range("H122").Value = 3012.25
$ws.Range("I122").Value = 2231.25
$ws.Range("J122").Value = 4574.25
$ws.Range("K122").Value = 6693.75
$ws.Range("L122").Value = 13722.75
$ws.Range("M122").Value = -4243.75
$ws.Range("N122").Value = -18622.75

# Row 132
$ws.Range("H132").Value = 903.5
$ws.Range("I132").Value = 992.75
$ws.Range("K132").Value = 2978.25
$ws.Range("M132").Value = -448.25

# Row 136
$ws.Range("H136").Value = 2818.8667
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()

# Row 38
$ws.Range("H38").Value = 6099
$ws.Range("I38").Value = 6099
$ws.Range("K38").Value = 6099
$ws.Range("M38").Value = -5683

# Row 86
$ws.Range("H86").Value = 2835.3
$ws.Range("I86").Value = 2692.3333
$ws.Range("J86").Value = 3049.75
$ws.Range("K86").Value = 2692.3333
$ws.Range("L86").Value = 3049.75
$ws.Range("M86").Value = -1569.3333
$ws.Range("N86").Value = -5295.75

# Row 89
$ws.Range("H89").Value = 2835.3
$ws.Range("I89").Value = 2692.3333
$ws.Range("J89").Value = 3049.75
$ws.Range("K89").Value = 13461.6665
$ws.Range("L89").Value = 15248.75
$ws.Range("M89").Value = -7845.666499999999
$ws.Range("N89").Value = -26480.75

# Row 105
$ws.Range("H105").Value = 4997.5
$ws.Range("I105").Value = 4997
$ws.Range("K105").Value = 4997
$ws.Range("M105").Value = -3250

# Row 106
$ws.Range("H106").Value = 50671
$ws.Range("J106").Value = 50671
$ws.Range("L106").Value = 50671
$ws.Range("N106").Value = -53195

# Row 107
$ws.Range("H107").Value = 3706.3333
$ws.Range("I107").Value = 2807.8235
$ws.Range("K107").Value = 2807.8235
$ws.Range("M107").Value = -887.8235

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5932.6665
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 31
$ws.Range("H31").Value = 2102.25
$ws.Range("I31").Value = 1660.7778
$ws.Range("K31").Value = 1660.7778
$ws.Range("M31").Value = -1365.7778

# Row 34
$ws.Range("H34").Value = 2102.25
$ws.Range("I34").Value = 1660.7778
$ws.Range("K34").Value = 1660.7778
$ws.Range("M34").Value = -1458.7778

# Row 62
$ws.Range("H62").Value = 4800
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4800
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6048

# Row 65
$ws.Range("H65").Value = 4800
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -30240

# Row 105
$ws.Range("H105").Value = 2797.5908
$ws.Range("I105").Value = 2302.8462
$ws.Range("K105").Value = 2302.8462
$ws.Range("M105").Value = -555.8462

# Row 113
$ws.Range("H113").Value = 5932.6665
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Row 132
$ws.Range("H132").Value = 2252.3809
$ws.Range("I132").Value = 2063.0557
$ws.Range("K132").Value = 6189.1671
$ws.Range("M132").Value = -3659.1671

$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 1451
$ws.Range("I41").Value = 1402
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 4206
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -3868
$ws.Range("N41").Value = -5176

# Row 70
$ws.Range("H70").Value = 9045.5
$ws.Range("I70").Value = 2159.6
$ws.Range("K70").Value = 6478.799999999999
$ws.Range("M70").Value = -6163.799999999999

# Row 73
$ws.Range("H73").Value = 9045.5
$ws.Range("I73").Value = 2159.6
$ws.Range("K73").Value = 6478.799999999999
$ws.Range("M73").Value = -5386.799999999999

# Row 121
$ws.Range("H121").Value = 6062
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 6062
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 18186
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -20806

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 125
$ws.Range("M2").Value = 112
$ws.Range("N2").Value = -351

# Row 126
$ws.Range("H126").Value = 4407
$ws.Range("I126").Value = 4423
$ws.Range("K126").Value = 13269
$ws.Range("M126").Value = -10799

# Row 132
$ws.Range("H132").Value = 1661
$ws.Range("I132").Value = 1604.6666
$ws.Range("K132").Value = 4813.9998
$ws.Range("M132").Value = -2283.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1103
$ws.Range("I16").Value = 1019.46155
$ws.Range("J16").Value = 1374.5
$ws.Range("K16").Value = 1019.46155
$ws.Range("L16").Value = 1374.5
$ws.Range("M16").Value = -849.46155
$ws.Range("N16").Value = -1714.5

# Row 46
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376

# Row 61
$ws.Range("H61").Value = 4624.625
$ws.Range("I61").Value = 4428.143
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 4428.143
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -4226.143
$ws.Range("N61").Value = -6404

# Row 93
$ws.Range("H93").Value = 1590.3334
$ws.Range("I93").Value = 1381.3077
$ws.Range("K93").Value = 1381.3077
$ws.Range("M93").Value = -133.3077000000001

# Row 113
$ws.Range("H113").Value = 4624.625
$ws.Range("I113").Value = 4428.143
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 4428.143
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -2258.143
$ws.Range("N113").Value = -10340

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 18837.8
$ws.Range("I62").Value = 17000.334
$ws.Range("K62").Value = 17000.334
$ws.Range("M62").Value = -16376.334

# Row 65
$ws.Range("H65").Value = 18837.8
$ws.Range("I65").Value = 17000.334
$ws.Range("K65").Value = 85001.67
$ws.Range("M65").Value = -81881.67

# Row 107
$ws.Range("H107").Value = 1541.4286
$ws.Range("I107").Value = 1155.4286
$ws.Range("J107").Value = 1927.4286
$ws.Range("K107").Value = 3466.2858
$ws.Range("L107").Value = 5782.2858
$ws.Range("M107").Value = -1546.2858
$ws.Range("N107").Value = -9622.2858

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 132
$ws.Range("H132").Value = 1695.0952
$ws.Range("I132").Value = 1735.4736
$ws.Range("K132").Value = 5206.4208
$ws.Range("M132").Value = -2676.4208

# Row 136
$ws.Range("H136").Value = 3560.8462
$ws.Range("I136").Value = 3580.5454
$ws.Range("J136").Value = 3452.5
$ws.Range("K136").Value = 10741.6362
$ws.Range("L136").Value = 10357.5
$ws.Range("M136").Value = -8191.636200000001
$ws.Range("N136").Value = -15457.5
